# Fix a typo in the "empresa" (company) paragraph of the agreement
# template: "${empresa.telefno}" -> "${empresa.telefono}".
#
# In the original document the placeholder text is split across runs as
#   ... "Teléfono ${empresa.t" | "ele" | <bookmark _GoBack/> | "fno}, e-mail ..."
# i.e. the word "telefno" is missing an "o" before "no}". The fix simply
# inserts the missing "o" right between the "f" and the "no}" that follow
# the existing "ele" run / "_GoBack" bookmark, turning the text into
# "${empresa.telefono}".

$d = $word.ActiveDocument

$inserted = $false

# Preferred approach: the last-edit bookmark "_GoBack" sits immediately
# after "...${empresa.tele" in the source document, i.e. right before the
# "fno}, e-mail ${empresa.email}. " run. Anchoring on it lets us place the
# new "o" precisely between the "f" and "no}" without disturbing anything
# else (and it naturally keeps the bookmark sitting right before the text
# that was just edited, exactly like Word leaves it after a real edit).
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks("_GoBack")
    $pos = $bm.Range.End + 1
    $target = $d.Range($pos, $pos)
    if ($target.Text -eq $null -and $true) {
        # no-op branch placeholder (collapsed ranges have no .Text issues,
        # kept for clarity of intent only)
    }
    $checkBefore = $d.Range($pos - 8, $pos)
    if ($checkBefore.Text -eq "empresa.tf" -or $checkBefore.Text -like "*tele*f") {
        $target.InsertBefore("o")
        $inserted = $true
    }
}

# Fallback: locate the text via Find (works even without the bookmark,
# and without merging/clobbering surrounding runs since no replacement
# text is supplied to Find.Execute -- it just moves/collapses the range).
if (-not $inserted) {
    $searchRange = $d.Content
    $found = $searchRange.Find.Execute("empresa.telef", $true, $false, $false, `
                                        $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $pos = $searchRange.End
        $target = $d.Range($pos, $pos)
        $target.InsertBefore("o")
        $inserted = $true
    }
}

if (-not $inserted) {
    throw "Could not locate '\${empresa.telefno}' to fix the typo."
}
